# Cyclically shift the data rows 3-7 (columns A:AY) down by one row,
# with row 7's data wrapping around to become the new row 3.
#
# Mapping (old -> new):
#   row 3 -> row 4
#   row 4 -> row 5
#   row 5 -> row 6
#   row 6 -> row 7
#   row 7 -> row 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 7
$lastCol = "AY"

# Columns that hold values which look numeric/date-like but must stay plain
# text (the workbook stores them as inline strings, e.g. "Antal" = "16" and
# "Startdatum"/"Slutdatum" = "2020-06-17"). Without forcing a text format,
# Excel's automatic type detection would turn them into real numbers/dates.
$textCols = @("I", "Y", "AA")

# Capture every source row's data (column A:AY) up front so that later writes
# never clobber data that still needs to be read.
$rowData = @{}
$rowTextVals = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData[$r] = $ws.Range("A" + $r + ":" + $lastCol + $r).Value2

    $colVals = @{}
    foreach ($col in $textCols) {
        $colVals[$col] = $ws.Range($col + $r).Value2
    }
    $rowTextVals[$r] = $colVals
}

# New row r gets the data that used to live in row (r-1), except the first
# row of the range, which wraps around and gets the data from the last row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -eq $firstRow) {
        $srcRow = $lastRow
    } else {
        $srcRow = $r - 1
    }

    $ws.Range("A" + $r + ":" + $lastCol + $r).Value2 = $rowData[$srcRow]

    # Re-apply the text-sensitive columns so their values remain text.
    foreach ($col in $textCols) {
        $cell = $ws.Range($col + $r)
        $cell.NumberFormat = "@"
        $cell.Value2 = $rowTextVals[$srcRow][$col]
        $cell.NumberFormat = "General"
    }
}
